$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Prezzi) values ---
$ws.Range("C2").Value = 115196.8512725545
$ws.Range("D2").Value = 510532.7437157534
$ws.Range("E2").Value = 28579.68403444847
$ws.Range("F2").Value = 539112.4277502019
$ws.Range("G2").Value = -31269.57633747661
$ws.Range("H2").Value = 507842.8514127253

# --- Update row 3 (Costi -> Costi totali) label and values ---
$ws.Range("A3").Value = "Costi totali"
$ws.Range("C3").Value = 85034.84442781005
$ws.Range("D3").Value = 376859.8877278101
$ws.Range("E3").Value = 36128.07294093538
$ws.Range("F3").Value = 412987.9606687455
$ws.Range("G3").Value = 74105.82933125459
$ws.Range("H3").Value = 487093.79

# --- Insert two new rows before the current MOL row (row 4), pushing it to row 6 ---
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Match formatting (border/bold/alignment style) of the label column from row 3
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- New row 4: Costi MP ---
$ws.Range("A4").Value = "Costi MP"
$ws.Range("B4").Value = 172660.4
$ws.Range("C4").Value = 50311.48145071982
$ws.Range("D4").Value = 222971.8814507198
$ws.Range("E4").Value = -10154.18096037768
$ws.Range("F4").Value = 212817.7004903422
$ws.Range("G4").Value = 46692.60950965786
$ws.Range("H4").Value = 259510.31

# --- New row 5: Costi risorse ---
$ws.Range("A5").Value = "Costi risorse"
$ws.Range("B5").Value = 119164.6433
$ws.Range("C5").Value = 34723.36297709029
$ws.Range("D5").Value = 153888.0062770903
$ws.Range("E5").Value = 46282.25390131297
$ws.Range("F5").Value = 200170.2601784033
$ws.Range("G5").Value = 27413.21982159675
$ws.Range("H5").Value = 227583.48

# --- Row 6 (was row 4): MOL, label stays the same, update values ---
$ws.Range("C6").Value = 30162.00684474444
$ws.Range("D6").Value = 133672.8559879434
$ws.Range("E6").Value = -7548.388906486914
$ws.Range("F6").Value = 126124.4670814564
$ws.Range("G6").Value = -105375.4056687312
$ws.Range("H6").Value = 20749.06141272525
